$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '68.394.53'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '3.921.50'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '486.23'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.39'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.624'
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000347'
$ws.Range('E11').Value = '  -2.95%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '43.19'
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.78'
$ws.Range('E13').Value = '  +2.90%  '
$ws.Range('D14').Value = '4.543.86'
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('D15').Value = '3.906.91'
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.35'
$ws.Range('E16').Value = '  -4.62%  '
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.11'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('D20').Value = '68.439.08'
$ws.Range('E20').Value = '  +1.10%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '433.85'
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '15.28'
$ws.Range('E22').Value = '  +4.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.50'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '88.40'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.61'
$ws.Range('E25').Value = '  +17.72%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.27'
$ws.Range('E26').Value = '  +11.73%  '
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '38.08'
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '715.60'
$ws.Range('E30').Value = '  -1.47%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.82'
$ws.Range('E31').Value = '  +2.89%  '
$ws.Range('E32').Value = '  -1.03%  '
$ws.Range('E33').Value = '  +4.09%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.16'
$ws.Range('E34').Value = '  +13.96%  '
$ws.Range('D35').Value = '0.0₃0886'
$ws.Range('E35').Value = '  +2.54%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '41.62'
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '60.95'
$ws.Range('E37').Value = '  +3.39%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.145'
$ws.Range('E39').Value = '  -5.37%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.395'
$ws.Range('E40').Value = '  +16.69%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.01'
$ws.Range('E41').Value = '  +18.88%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0495'
$ws.Range('E42').Value = '  +3.87%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.12'
$ws.Range('E43').Value = '  +2.47%  '
$ws.Range('E44').Value = '  +4.61%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.40'
$ws.Range('E45').Value = '  +5.31%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.143'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.15'
$ws.Range('E49').Value = '  -5.42%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0348'
$ws.Range('E50').Value = '  +27.85%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '145.12'
$ws.Range('E51').Value = '  -2.48%  '
